$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full NATMI Ngf -> Sort1 table, recomputed to include the "ECs" sending/target
# cluster (per Dr Hou's advice). Replaces the previous 10-row (FAPs/sCs) table
# with a 15-row table spanning ECs / FAPs / sCs senders x ECs / FAPs / M1 / M2 / sCs
# targets.

$data = @(
        @("ECs","Ngf","Sort1","ECs",1,0.3333333333333333,0.672676,2.018028,0.1010994183758018,0.1010994183758018,3,1,0.72097,2.16291,0.03318293102712134,0.03318293102712135,0.48497921572,4.36481294148,0.003354775026846314,0.003354775026846315),
        @("ECs","Ngf","Sort1","FAPs",1,0.3333333333333333,0.672676,2.018028,0.1010994183758018,0.1010994183758018,3,1,1.851391,5.554173,0.08521100719479757,0.08521100719479759,1.245386292316,11.208476630844,0.008614783266610295,0.008614783266610297),
        @("ECs","Ngf","Sort1","M1",1,0.3333333333333333,0.672676,2.018028,0.1010994183758018,0.1010994183758018,3,1,8.824443666666665,26.473331,0.4061485298191571,0.4061485298191572,5.935991467918666,53.423923211268,0.04106138013890376,0.04106138013890378),
        @("ECs","Ngf","Sort1","M2",1,0.3333333333333333,0.672676,2.018028,0.1010994183758018,0.1010994183758018,3,1,8.620174666666665,25.860524,0.3967469678429597,0.3967469678429598,5.798584614074667,52.187261526672,0.04011088769128616,0.04011088769128617),
        @("ECs","Ngf","Sort1","sCs",1,0.3333333333333333,0.672676,2.018028,0.1010994183758018,0.1010994183758018,3,1,1.710155,5.130465,0.0787105641159642,0.0787105641159642,1.15038022478,10.35342202302,0.007957592252155234,0.007957592252155234),
        @("FAPs","Ngf","Sort1","ECs",3,1,3.536689,10.610067,0.5315444595557087,0.5315444595557088,3,1,0.72097,2.16291,0.03318293102712134,0.03318293102712135,2.54984666833,22.94862001497,0.01763820313928557,0.01763820313928558),
        @("FAPs","Ngf","Sort1","FAPs",3,1,3.536689,10.610067,0.5315444595557087,0.5315444595557088,3,1,1.851391,5.554173,0.08521100719479757,0.08521100719479759,6.547794184398998,58.93014765959099,0.04529343876755628,0.0452934387675563),
        @("FAPs","Ngf","Sort1","M1",3,1,3.536689,10.610067,0.5315444595557087,0.5315444595557088,3,1,8.824443666666665,26.473331,0.4061485298191571,0.4061485298191572,31.20931284701966,280.8838156231769,0.2158860007820695,0.2158860007820696),
        @("FAPs","Ngf","Sort1","M2",3,1,3.536689,10.610067,0.5315444595557087,0.5315444595557088,3,1,8.620174666666665,25.860524,0.3967469678429597,0.3967469678429598,30.48687692167866,274.3818922951079,0.2108886526024522,0.2108886526024522),
        @("FAPs","Ngf","Sort1","sCs",3,1,3.536689,10.610067,0.5315444595557087,0.5315444595557088,3,1,1.710155,5.130465,0.0787105641159642,0.0787105641159642,6.048286376795,54.434577391155,0.04183816426434515,0.04183816426434515),
        @("sCs","Ngf","Sort1","ECs",3,1,2.444244,7.332732,0.3673561220684894,0.3673561220684895,3,1,0.72097,2.16291,0.03318293102712134,0.03318293102712135,1.76222659668,15.86003937012,0.01218995286098945,0.01218995286098946),
        @("sCs","Ngf","Sort1","FAPs",3,1,2.444244,7.332732,0.3673561220684894,0.3673561220684895,3,1,1.851391,5.554173,0.08521100719479757,0.08521100719479759,4.525251343403999,40.727262090636,0.03130278516063099,0.031302785160631),
        @("sCs","Ngf","Sort1","M1",3,1,2.444244,7.332732,0.3673561220684894,0.3673561220684895,3,1,8.824443666666665,26.473331,0.4061485298191571,0.4061485298191572,21.569093485588,194.121841370292,0.1492011488981838,0.1492011488981839),
        @("sCs","Ngf","Sort1","M2",3,1,2.444244,7.332732,0.3673561220684894,0.3673561220684895,3,1,8.620174666666665,25.860524,0.3967469678429597,0.3967469678429598,21.069810207952,189.628291871568,0.1457474275492214,0.1457474275492214),
        @("sCs","Ngf","Sort1","sCs",3,1,2.444244,7.332732,0.3673561220684894,0.3673561220684895,3,1,1.710155,5.130465,0.0787105641159642,0.0787105641159642,4.18003609782,37.62032488038,0.02891480759946381,0.02891480759946381)
)

$startRow = 2
$numRows = $data.Length
$numCols = 20

# Clear out any previous contents in the body so a shrink/grow in row count
# never leaves stale cells behind, then write the full replacement block.
$ws.Range("A2:T100").ClearContents()

for ($i = 0; $i -lt $numRows; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($j = 0; $j -lt $numCols; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
